$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=21.24103426297825; "C"=10.16730127147443; "D"=6.783034899723157; "E"=9.242673668494744; "F"=36.50089427488133; "H"=7.344005520526261; "I"=28.57412395558024; "L"=10.51400227631766; "N"=18.8490871953992 }
    3 = @{ "B"=20.7331130343156; "C"=9.571374779614194; "D"=6.813626447635657; "E"=9.258936321943818; "F"=36.25053531943155; "H"=7.344005520526261; "I"=28.5699713761728; "L"=10.49248674702796; "N"=18.92323688013329 }
    4 = @{ "B"=20.42113794443141; "C"=9.188186371296117; "D"=6.833358383626037; "E"=9.26955483221135; "F"=36.10809446377849; "H"=7.344005520526261; "I"=28.57515571999571; "L"=10.48153583486208; "N"=18.97075048805318 }
    5 = @{ "B"=20.29418125165343; "C"=9.027830497371996; "D"=6.841638131436982; "E"=9.274041532107173; "F"=36.05292966471112; "H"=7.344005520526261; "I"=28.57920677570069; "L"=10.47764380748857; "N"=18.99061338602415 }
    6 = @{ "B"=20.27311685253524; "C"=9.000954446605586; "D"=6.843027410418425; "E"=9.274796193342731; "F"=36.04394478096481; "H"=7.344005520526261; "I"=28.57999631239072; "L"=10.47703205817565; "N"=18.99394189045201 }
    7 = @{ "B"=20.4194247779991; "C"=9.186040563501903; "D"=6.833469079852981; "E"=9.269614694835154; "F"=36.1073387728125; "H"=7.344005520526261; "I"=28.57520251507655; "L"=10.48148103274049; "N"=18.9710163367843 }
    8 = @{ "B"=21.06604568684344; "C"=9.965490839614146; "D"=6.793386155724206; "E"=9.248149828323216; "F"=36.41225763017709; "H"=7.344005520526261; "I"=28.57108347763698; "L"=10.50611626462322; "N"=18.87424282729003 }
    9 = @{ "B"=22.32462432897404; "C"=11.35190367399198; "D"=6.722297080481603; "E"=9.211065352732946; "F"=37.09748599449368; "H"=7.344005520526261; "I"=28.62460295030326; "L"=10.57223251628998; "N"=18.70015472599574 }
    10 = @{ "B"=23.23267697735591; "C"=12.29779322933998; "D"=6.674629660981654; "E"=9.186850531309052; "F"=37.65079809015274; "H"=7.344005520526261; "I"=28.70170920611979; "L"=10.631469094186; "N"=18.5817152093954 }
    11 = @{ "B"=23.64006818330384; "C"=12.72989003967784; "D"=6.653931447733713; "E"=9.176488174596102; "F"=37.91253805314977; "H"=7.344005520526261; "I"=28.74501282236967; "L"=10.66068086674379; "N"=18.5298677778927 }
    12 = @{ "B"=23.79335975593877; "C"=12.88978341918996; "D"=6.646235152498687; "E"=9.17265778845773; "F"=38.0130218233759; "H"=7.344005520526261; "I"=28.76259308836812; "L"=10.67206313834741; "N"=18.51052507747634 }
    13 = @{ "B"=23.76039163344799; "C"=12.85551321178842; "D"=6.647886388050159; "E"=9.173478571885084; "F"=37.99132119941108; "H"=7.344005520526261; "I"=28.75875431951466; "L"=10.66959760086713; "N"=18.51467796135181 }
    14 = @{ "B"=23.6527000658852; "C"=12.74311925640055; "D"=6.653295430164585; "E"=9.176171171967603; "F"=37.92077792404125; "H"=7.344005520526261; "I"=28.74643548077372; "L"=10.6616109033785; "N"=18.52827062168794 }
    15 = @{ "B"=23.58660385160453; "C"=12.67378922067831; "D"=6.656627071353826; "E"=9.177832649696176; "F"=37.87774408802373; "H"=7.344005520526261; "I"=28.73904372955514; "L"=10.65676038318135; "N"=18.53663435059098 }
    16 = @{ "B"=23.20592454904557; "C"=12.26903137171112; "D"=6.676002167679412; "E"=9.187540851740049; "F"=37.63388840599678; "H"=7.344005520526261; "I"=28.69904472196442; "L"=10.6296051694093; "N"=18.58514432422558 }
    17 = @{ "B"=22.9708154097782; "C"=12.01653678993965; "D"=6.688140599012252; "E"=9.193663575672575; "F"=37.48680914821198; "H"=7.344005520526261; "I"=28.6766136714961; "L"=10.61352278486977; "N"=18.6154229543176 }
    18 = @{ "B"=22.8350604864836; "C"=11.87903912708889; "D"=6.695215109444193; "E"=9.197246693790573; "F"=37.40316232375717; "H"=7.344005520526261; "I"=28.66448649774829; "L"=10.60448613667654; "N"=18.63302969374395 }
    19 = @{ "B"=22.78901089657314; "C"=11.83215378823488; "D"=6.69762636025643; "E"=9.198470446069065; "F"=37.37500620465543; "H"=7.344005520526261; "I"=28.66051348478771; "L"=10.60146330695472; "N"=18.63902392427745 }
    20 = @{ "B"=22.99589885933153; "C"=12.04182693240931; "D"=6.686838838091015; "E"=9.193005439541743; "F"=37.50236827021786; "H"=7.344005520526261; "I"=28.67892133226508; "L"=10.61521272094731; "N"=18.61217995599825 }
    21 = @{ "B"=23.68435947295149; "C"=12.77623322788247; "D"=6.651702819253306; "E"=9.175377751775006; "F"=37.94146166194646; "H"=7.344005520526261; "I"=28.75002175493069; "L"=10.66394813843517; "N"=18.52427024599892 }
    22 = @{ "B"=24.1285378396831; "C"=13.23470802046167; "D"=6.629565122831584; "E"=9.164402554231556; "F"=38.23637478608247; "H"=7.344005520526261; "I"=28.8033788326155; "L"=10.69766437610733; "N"=18.46851056266725 }
    23 = @{ "B"=23.89204955847352; "C"=12.99199486383812; "D"=6.64130490319822; "E"=9.170210406697921; "F"=38.07827296492234; "H"=7.344005520526261; "I"=28.77427153328335; "L"=10.67950060750461; "N"=18.49811594326396 }
    24 = @{ "B"=22.98456044952684; "C"=12.03039947664983; "D"=6.687427065286798; "E"=9.193302786375703; "F"=37.49533114797817; "H"=7.344005520526261; "I"=28.67787564433211; "L"=10.61444804844027; "N"=18.61364549512675 }
    25 = @{ "B"=21.98631678544683; "C"=10.99264676656235; "D"=6.740726206918011; "E"=9.220563821023388; "F"=36.9030945633605; "H"=7.344005520526261; "I"=28.60350039645846; "L"=10.55245811483969; "N"=18.51467796135181 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}